$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old table area (A1:D3) before writing the new layout
$ws.Range("A1:D3").Clear()

# New header row
$ws.Range("A1").Value = "posFile"
$ws.Range("B1").Value = "stimFile"
$ws.Range("C1").Value = "cs_plus_s"
$ws.Range("D1").Value = "cs_minus_s"
$ws.Range("E1").Value = "cs_plus_ns"
$ws.Range("F1").Value = "cs_minus_ns"

# New data row
$ws.Range("A2").Value = "positions.xlsx"
$ws.Range("B2").Value = "stimuli.xlsx"
$ws.Range("C2").Value = "stimuli/social/016_y_m_n_b.jpg"
$ws.Range("D2").Value = "stimuli/social/031_y_m_n_a.jpg"
$ws.Range("E2").Value = "stimuli/non-social/016_y_m_n_b_scrambled.jpg"
$ws.Range("F2").Value = "stimuli/non-social/031_y_m_n_a_scrambled.jpg"

# Adjust column widths to match the new content layout (closest reachable
# values given the runtime's internal pixel-rounding of ColumnWidth)
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666
$ws.Columns.Item(3).ColumnWidth = 40.0
$ws.Columns.Item(4).ColumnWidth = 26.333333333333332
$ws.Columns.Item(5).ColumnWidth = 40.0
$ws.Columns.Item(6).ColumnWidth = 39.833333333333336

$ws.Range("E11").Select() | Out-Null
